$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("J6").Value = "."
$ws.Range("J7").Select()
